$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.951.61'
$ws.Range('E2').Value = '  -1.02%  '
$ws.Range('D3').Value = '2.928.59'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '359.79'
$ws.Range('E5').Value = '  +2.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '110.40'
$ws.Range('E6').Value = '  -2.45%  '
$ws.Range('E7').Value = '  +1.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.632'
$ws.Range('E9').Value = '  +0.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.42'
$ws.Range('E10').Value = '  -2.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0880'
$ws.Range('E12').Value = '  +0.91%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.66'
$ws.Range('E13').Value = '  -2.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.89'
$ws.Range('E14').Value = '  -0.07%  '
$ws.Range('D15').Value = '3.390.97'
$ws.Range('E15').Value = '  +0.08%  '
$ws.Range('D16').Value = '2.938.95'
$ws.Range('E16').Value = '  +0.26%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.988'
$ws.Range('E17').Value = '  -0.48%  '
$ws.Range('D18').Value = '51.977.75'
$ws.Range('E18').Value = '  -1.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.37'
$ws.Range('E19').Value = '  +0.90%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.62'
$ws.Range('E20').Value = '  -1.35%  '
$ws.Range('E21').Value = '  -2.80%  '
$ws.Range('D22').Value = '0.0₃0985'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '71.15'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '270.26'
$ws.Range('E24').Value = '  -0.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.85'
$ws.Range('E25').Value = '  +1.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.187'
$ws.Range('E26').Value = '  +13.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '27.11'
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.60'
$ws.Range('E28').Value = '  +16.18%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.108'
$ws.Range('E30').Value = '  +14.42%  '
$ws.Range('E31').Value = '  -0.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '38.49'
$ws.Range('E32').Value = '  +1.02%  '
$ws.Range('E33').Value = '  +1.65%  '
$ws.Range('E34').Value = '  -1.87%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '52.14'
$ws.Range('E35').Value = '  -1.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0446'
$ws.Range('E36').Value = '  -1.78%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E38').Value = '  -2.87%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.47'
$ws.Range('E39').Value = '  -2.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.02'
$ws.Range('E40').Value = '  -3.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.74'
$ws.Range('E41').Value = '  -0.27%  '
$ws.Range('E42').Value = '  +2.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '23.07'
$ws.Range('E43').Value = '  -6.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '119.47'
$ws.Range('E44').Value = '  -2.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.17'
$ws.Range('E45').Value = '  -1.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.50'
$ws.Range('E46').Value = '  -2.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.48'
$ws.Range('E47').Value = '  -5.06%  '
$ws.Range('D48').Value = '2.134.87'
$ws.Range('E48').Value = '  -3.81%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.251'
$ws.Range('E49').Value = '  -5.55%  '
$ws.Range('E50').Value = '  -1.78%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.19'
$ws.Range('E51').Value = '  +0.02%  '
